$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Quantity" column before the existing "Outcome Status" column (column I)
$ws.Columns.Item(9).Insert()

# Set the new column's header
$ws.Range("I1").Value = "Quantity"

# Populate the Quantity values for the existing appointment rows
$ws.Range("I2").Value = "2,1"
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 0

# Remove the duplicate/stale appointment record (previously row 5)
$ws.Rows.Item(5).Delete()

# Restore the active cell selection
$null = $ws.Range("H9").Select()
